# Implemented getting kafka relations.
# The classFields extraction now walks fields in a different (declaration)
# order for a few classes, which re-shuffles the field rows for
# OrderGeneratorService, OrderControllerTests and OrderController on the
# "classFields" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("classFields")

# pl.piomin.order.service.OrderGeneratorService (rows 4-7):
# was: template, executor, RAND, id
# now: executor, RAND, id, template
$ws.Cells.Item(4, 2).Value2 = "executor"
$ws.Cells.Item(4, 3).Value2 = "private"
$ws.Cells.Item(4, 4).Value2 = "java.util.concurrent.Executor"

$ws.Cells.Item(5, 2).Value2 = "RAND"
$ws.Cells.Item(5, 3).Value2 = "private"
$ws.Cells.Item(5, 4).Value2 = "java.util.Random"

$ws.Cells.Item(6, 2).Value2 = "id"
$ws.Cells.Item(6, 3).Value2 = "private"
$ws.Cells.Item(6, 4).Value2 = "java.util.concurrent.atomic.AtomicLong"

$ws.Cells.Item(7, 2).Value2 = "template"
$ws.Cells.Item(7, 3).Value2 = "private"
$ws.Cells.Item(7, 4).Value2 = "org.springframework.kafka.core.KafkaTemplate"

# pl.piomin.order.OrderControllerTests (rows 9-10):
# was: mapper, template
# now: template, mapper
$ws.Cells.Item(9, 2).Value2 = "template"
$ws.Cells.Item(9, 3).Value2 = "private"
$ws.Cells.Item(9, 4).Value2 = "org.springframework.kafka.core.KafkaTemplate"

$ws.Cells.Item(10, 2).Value2 = "mapper"
$ws.Cells.Item(10, 3).Value2 = ""
$ws.Cells.Item(10, 4).Value2 = "com.fasterxml.jackson.databind.ObjectMapper"

# pl.piomin.order.controller.OrderController (rows 12-15):
# was: LOG, id, template, orderGeneratorService
# now: orderGeneratorService, LOG, id, template
$ws.Cells.Item(12, 2).Value2 = "orderGeneratorService"
$ws.Cells.Item(12, 3).Value2 = "private"
$ws.Cells.Item(12, 4).Value2 = "pl.piomin.order.service.OrderGeneratorService"

$ws.Cells.Item(13, 2).Value2 = "LOG"
$ws.Cells.Item(13, 3).Value2 = "private"
$ws.Cells.Item(13, 4).Value2 = "org.slf4j.Logger"

$ws.Cells.Item(14, 2).Value2 = "id"
$ws.Cells.Item(14, 3).Value2 = "private"
$ws.Cells.Item(14, 4).Value2 = "java.util.concurrent.atomic.AtomicLong"

$ws.Cells.Item(15, 2).Value2 = "template"
$ws.Cells.Item(15, 3).Value2 = "private"
$ws.Cells.Item(15, 4).Value2 = "org.springframework.kafka.core.KafkaTemplate"
